$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '23.713.07'
$ws.Range("E2").Value = '  +1.28%  '

# Row 3
$ws.Range("D3").Value = '1.657.08'
$ws.Range("E3").Value = '  +1.10%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.20%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  -0.12%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.18'
$ws.Range("E6").Value = '  -0.11%  '

# Row 7
$ws.Range("E7").Value = '  +0.50%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3610'
$ws.Range("E8").Value = '  +0.39%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.08'
$ws.Range("E9").Value = '  -1.58%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08201'
$ws.Range("E10").Value = '  +0.02%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.232'
$ws.Range("E11").Value = '  +0.37%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.22%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.42'
$ws.Range("E13").Value = '  -0.18%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.455'
$ws.Range("E14").Value = '  +0.46%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.448'
$ws.Range("E15").Value = '  +1.68%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001222'
$ws.Range("E16").Value = '  -1.07%  '

# Row 17
$ws.Range("D17").Value = '1.656.83'
$ws.Range("E17").Value = '  +1.34%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.69'
$ws.Range("E18").Value = '  +2.65%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07037'
$ws.Range("E19").Value = '  +1.02%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.790'
$ws.Range("E20").Value = '  +3.17%  '

# Row 21
$ws.Range("E21").Value = '  +0.77%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.08%  '

# Row 23
$ws.Range("E23").Value = '  +1.71%  '

# Row 24
$ws.Range("D24").Value = '23.716.66'
$ws.Range("E24").Value = '  +1.32%  '

# Row 25
$ws.Range("E25").Value = '  -2.17%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.029'
$ws.Range("E26").Value = '  -0.79%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.27'
$ws.Range("E27").Value = '  +0.63%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.49'
$ws.Range("E28").Value = '  +1.22%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.242'
$ws.Range("E29").Value = '  -0.55%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.92'
$ws.Range("E30").Value = '  -0.26%  '

# Row 31
$ws.Range("D31").Value = '1.840.30'
$ws.Range("E31").Value = '  +1.21%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.102'
$ws.Range("E32").Value = '  +9.45%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.254'
$ws.Range("E33").Value = '  +4.89%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.00'
$ws.Range("E34").Value = '  +5.05%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.057'
$ws.Range("E35").Value = '  -2.84%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02809'

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2505'
$ws.Range("E37").Value = '  +0.06%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08815'
$ws.Range("E38").Value = '  +0.54%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.086'
$ws.Range("E39").Value = '  +2.33%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06991'
$ws.Range("E40").Value = '  -0.39%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.04'
$ws.Range("E41").Value = '  +6.84%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6992'
$ws.Range("E42").Value = '  -0.25%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.338'
$ws.Range("E43").Value = '  -0.59%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.98'
$ws.Range("E44").Value = '  +2.56%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6510'
$ws.Range("E45").Value = '  +0.34%  '

# Row 46
$ws.Range("E46").Value = '  -0.07%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.298'
$ws.Range("E47").Value = '  +0.60%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.966'
$ws.Range("E48").Value = '  +0.14%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07901'
$ws.Range("E49").Value = '  -0.76%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.26'
$ws.Range("E50").Value = '  -0.19%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.179'
$ws.Range("E51").Value = '  -0.97%  '

